$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update quantities for Female DB25 (row 6) and Male DB25 (row 7) from 4 to 3
$ws.Range("D6").Value = 3
$ws.Range("D7").Value = 3

# Update the active selection to D9
$ws.Range("D9").Select()

$excel.Calculate()
